$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Summary statistics formulas (E12:E16) referencing the table
$ws.Range("E12").Formula = "=SUM(Table1[Total Selling Price])"
$ws.Range("E13").Formula = "=AVERAGE(Table1[Total Selling Price])"
$ws.Range("E14").Formula = "=MAX(Table1[Profit])"
$ws.Range("E15").Formula = "=MIN(Table1[Profit])"
$ws.Range("E16").Formula = "=_xlfn.STDEV.S(Table1[Profit])"

# Labels for highest / lowest profit rows (column H)
$ws.Range("H3").Value = "highest profit"
$ws.Range("H5").Value = "lowest profit"

# Column H width
$ws.Columns.Item(8).ColumnWidth = 14.625

# Match the style used on H3/H5 (fill color) to the existing "D12/D13" style family
$ws.Range("H3:H5").Interior.ThemeColor = 7
$ws.Range("H3:H5").Interior.TintAndShade = 0.79998168889431442

# Move the active selection to E16
$ws.Range("E16").Select()
